# Auto-generated edit script: updates cryptocurrency price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''24.651.27'
$ws.Range("E2").Value = '  +1.82%  '

$ws.Range("D3").Value = '''1.702.58'
$ws.Range("E3").Value = '  +1.66%  '

$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''309.25'
$ws.Range("E5").Value = '  +0.01%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").Value = '''0.3731'
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("D8").Value = '''48.89'
$ws.Range("E8").Value = '  +2.37%  '

$ws.Range("D9").Value = '''0.3416'
$ws.Range("E9").Value = '  -0.94%  '

$ws.Range("D10").Value = '''1.177'
$ws.Range("E10").Value = '  -0.31%  '

$ws.Range("D11").Value = '''0.07427'
$ws.Range("E11").Value = '  +1.61%  '

$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  +0.24%  '

$ws.Range("D13").Value = '''20.77'
$ws.Range("E13").Value = '  +1.70%  '

$ws.Range("D14").Value = '''6.212'
$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").Value = '''6.888'
$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("D16").Value = '''1.694.85'
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").Value = '''0.00001115'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").Value = '''1.001'
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").Value = '''0.06681'
$ws.Range("E19").Value = '  -0.74%  '

$ws.Range("D20").Value = '''82.92'
$ws.Range("E20").Value = '  +1.30%  '

$ws.Range("D21").Value = '''17.02'
$ws.Range("E21").Value = '  +3.16%  '

$ws.Range("D22").Value = '''6.315'
$ws.Range("E22").Value = '  +2.69%  '

$ws.Range("D23").Value = '''13.05'
$ws.Range("E23").Value = '  +8.62%  '

$ws.Range("D24").Value = '''24.591.56'
$ws.Range("E24").Value = '  +1.76%  '

$ws.Range("D25").Value = '''2.437'
$ws.Range("E25").Value = '  +1.22%  '

$ws.Range("D26").Value = '''2.754'
$ws.Range("E26").Value = '  +3.55%  '

$ws.Range("D27").Value = '''20.10'
$ws.Range("E27").Value = '  +2.86%  '

$ws.Range("D28").Value = '''149.52'
$ws.Range("E28").Value = '  -1.41%  '

$ws.Range("D29").Value = '''130.79'
$ws.Range("E29").Value = '  +3.00%  '

$ws.Range("D30").Value = '''1.883.41'
$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("D31").Value = '''1.166'
$ws.Range("E31").Value = '  +17.60%  '

$ws.Range("D32").Value = '''6.644'
$ws.Range("E32").Value = '  +3.20%  '

$ws.Range("D33").Value = '''4.212'
$ws.Range("E33").Value = '  +2.37%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '''0.08727'
$ws.Range("E34").Value = '  +2.57%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '''1.761'
$ws.Range("E35").Value = '  -0.64%  '

$ws.Range("D36").Value = '''13.52'
$ws.Range("E36").Value = '  +7.37%  '

$ws.Range("D37").Value = '''5.460'
$ws.Range("E37").Value = '  +1.69%  '

$ws.Range("D38").Value = '''0.06468'
$ws.Range("E38").Value = '  -0.31%  '

$ws.Range("D39").Value = '''0.02364'
$ws.Range("E39").Value = '  +0.15%  '

$ws.Range("D40").Value = '''8.869'
$ws.Range("E40").Value = '  -0.37%  '

$ws.Range("D41").Value = '''0.2178'
$ws.Range("E41").Value = '  +1.78%  '

$ws.Range("D42").Value = '''1.271'
$ws.Range("E42").Value = '  -0.65%  '

$ws.Range("D43").Value = '''0.6381'
$ws.Range("E43").Value = '  +2.95%  '

$ws.Range("D44").Value = '''1.002'
$ws.Range("E44").Value = '  +0.31%  '

$ws.Range("D45").Value = '''13.74'
$ws.Range("E45").Value = '  +3.94%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '''0.6043'
$ws.Range("E46").Value = '  +1.64%  '

$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = '''3.805'
$ws.Range("E47").Value = '  -0.11%  '

$ws.Range("D48").Value = '''2.101'
$ws.Range("E48").Value = '  +3.40%  '

$ws.Range("D49").Value = '''128.57'
$ws.Range("E49").Value = '  +0.96%  '

$ws.Range("D50").Value = '''0.07210'
$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").Value = '''78.57'
$ws.Range("E51").Value = '  +2.54%  '
